# Update 杭州-漫展信息 workbook: refresh "想去人数" (want-to-go count) figures
# and flip one "已售罄" (sold out) status to "不可售" (not available) across
# the 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2842
$ws1.Range("F3").Value = 1591
$ws1.Range("F4").Value = 1162
$ws1.Range("G5").Value = "不可售"
$ws1.Range("F6").Value = 9690
$ws1.Range("F13").Value = 705
$ws1.Range("F15").Value = 1210
$ws1.Range("F16").Value = 1008
$ws1.Range("F17").Value = 3000
$ws1.Range("F18").Value = 2272
$ws1.Range("F20").Value = 1973
$ws1.Range("F22").Value = 1937
$ws1.Range("F24").Value = 1566
$ws1.Range("F25").Value = 313
$ws1.Range("F26").Value = 24
$ws1.Range("F27").Value = 195
$ws1.Range("F28").Value = 223
$ws1.Range("F30").Value = 347
$ws1.Range("F32").Value = 321
$ws1.Range("F33").Value = 527
$ws1.Range("F34").Value = 37
$ws1.Range("F35").Value = 160
$ws1.Range("F36").Value = 1550
$ws1.Range("F37").Value = 196
$ws1.Range("F38").Value = 1531
$ws1.Range("F39").Value = 55
$ws1.Range("F40").Value = 363
$ws1.Range("F41").Value = 31
$ws1.Range("F42").Value = 386
$ws1.Range("F43").Value = 798

# --- 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 29

# --- 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2842
$ws4.Range("F3").Value = 1591
$ws4.Range("G4").Value = "不可售"
$ws4.Range("F5").Value = 9690
$ws4.Range("F14").Value = 705
$ws4.Range("F15").Value = 1210
$ws4.Range("F16").Value = 1008
$ws4.Range("F17").Value = 3000
$ws4.Range("F18").Value = 2272
$ws4.Range("F19").Value = 1973
$ws4.Range("F20").Value = 1937
$ws4.Range("F22").Value = 1566
$ws4.Range("F23").Value = 313
$ws4.Range("F24").Value = 24
$ws4.Range("F25").Value = 195
$ws4.Range("F26").Value = 223
$ws4.Range("F28").Value = 347
$ws4.Range("F30").Value = 321
$ws4.Range("F31").Value = 527
$ws4.Range("F32").Value = 29
$ws4.Range("F35").Value = 37
$ws4.Range("F36").Value = 160
$ws4.Range("F37").Value = 1550
$ws4.Range("F39").Value = 196
$ws4.Range("F40").Value = 1531
$ws4.Range("F41").Value = 55
$ws4.Range("F43").Value = 363
$ws4.Range("F44").Value = 31
$ws4.Range("F45").Value = 386
$ws4.Range("F46").Value = 798

